$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 345.8889
$ws.Range("I53").Value = 262.66666
$ws.Range("J53").Value = 512.3333
$ws.Range("K53").Value = 262.66666
$ws.Range("L53").Value = 512.3333
$ws.Range("M53").Value = 374.33334
$ws.Range("N53").Value = -1786.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 284.77777
$ws.Range("I55").Value = 314.75
$ws.Range("J55").Value = 260.8
$ws.Range("K55").Value = 314.75
$ws.Range("L55").Value = 260.8
$ws.Range("M55").Value = -100.75
$ws.Range("N55").Value = -688.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2240.7273
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 2364.8
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 2364.8
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -3176.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2240.7273
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 2364.8
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 2364.8
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -5172.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3434.65
$ws.Range("I132").Value = 3135.2144
$ws.Range("J132").Value = 4133.3335
$ws.Range("K132").Value = 9405.643199999999
$ws.Range("L132").Value = 12400.0005
$ws.Range("M132").Value = -6875.643199999999
$ws.Range("N132").Value = -17460.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4642.5713
$ws.Range("I137").Value = 4499.5
$ws.Range("J137").Value = 4833.3335
$ws.Range("K137").Value = 13498.5
$ws.Range("L137").Value = 14500.0005
$ws.Range("M137").Value = -10948.5
$ws.Range("N137").Value = -19600.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 3500
$ws.Range("I17").Value = 3000
$ws.Range("J17").Value = 4000
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 4000
$ws.Range("M17").Value = -2827
$ws.Range("N17").Value = -4346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3187.077
$ws.Range("I45").Value = 1062.1428
$ws.Range("K45").Value = 1062.1428
$ws.Range("M45").Value = -685.1428000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 30249.6
$ws.Range("J62").Value = 30249.6
$ws.Range("L62").Value = 30249.6
$ws.Range("N62").Value = -31497.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6356.4287
$ws.Range("I63").Value = 2747.5
$ws.Range("J63").Value = 7800
$ws.Range("K63").Value = 2747.5
$ws.Range("L63").Value = 7800
$ws.Range("M63").Value = -2061.5
$ws.Range("N63").Value = -9172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 30249.6
$ws.Range("J65").Value = 30249.6
$ws.Range("L65").Value = 90748.79999999999
$ws.Range("N65").Value = -96988.79999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6356.4287
$ws.Range("I66").Value = 2747.5
$ws.Range("J66").Value = 7800
$ws.Range("K66").Value = 13737.5
$ws.Range("L66").Value = 39000
$ws.Range("M66").Value = -10305.5
$ws.Range("N66").Value = -45864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1987
$ws.Range("I132").Value = 1987
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5961
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3431
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 25019000
$ws.Range("J135").Value = 38000
$ws.Range("L135").Value = 38000
$ws.Range("N135").Value = -48140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 229
$ws.Range("I16").Value = 229
$ws.Range("K16").Value = 229
$ws.Range("M16").Value = -59

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 5000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -5346

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4354.304
$ws.Range("I86").Value = 3694.0588
$ws.Range("J86").Value = 6225
$ws.Range("K86").Value = 3694.0588
$ws.Range("L86").Value = 6225
$ws.Range("M86").Value = -2571.0588
$ws.Range("N86").Value = -8471

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4354.304
$ws.Range("I89").Value = 3694.0588
$ws.Range("J89").Value = 6225
$ws.Range("K89").Value = 18470.294
$ws.Range("L89").Value = 31125
$ws.Range("M89").Value = -12854.294
$ws.Range("N89").Value = -42357

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2410.4443
$ws.Range("I134").Value = 2198.75
$ws.Range("J134").Value = 2579.8
$ws.Range("K134").Value = 6596.25
$ws.Range("L134").Value = 7739.400000000001
$ws.Range("M134").Value = -4061.25
$ws.Range("N134").Value = -12809.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 23121494
$ws.Range("J86").Value = 7999
$ws.Range("L86").Value = 7999
$ws.Range("N86").Value = -10245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 23121494
$ws.Range("J89").Value = 7999
$ws.Range("L89").Value = 39995
$ws.Range("N89").Value = -51227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1846
$ws.Range("I134").Value = 1276.25
$ws.Range("J134").Value = 2301.8
$ws.Range("K134").Value = 3828.75
$ws.Range("L134").Value = 6905.400000000001
$ws.Range("M134").Value = -1293.75
$ws.Range("N134").Value = -11975.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 214.33333
$ws.Range("J15").Value = 239
$ws.Range("L15").Value = 717
$ws.Range("N15").Value = -997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 200
$ws.Range("J16").Value = 200
$ws.Range("L16").Value = 600
$ws.Range("N16").Value = -946

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1599
$ws.Range("J17").Value = 1599
$ws.Range("L17").Value = 4797
$ws.Range("N17").Value = -5135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 7505.25
$ws.Range("J19").Value = 7505.25
$ws.Range("L19").Value = 22515.75
$ws.Range("N19").Value = -22863.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 836.8
$ws.Range("J34").Value = 998.4
$ws.Range("L34").Value = 2995.2
$ws.Range("N34").Value = -3163.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("K39").Value = 1500
$ws.Range("M39").Value = -1206

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1917
$ws.Range("J46").Value = 2199.7334
$ws.Range("L46").Value = 2199.7334
$ws.Range("N46").Value = -2575.7334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 26250
$ws.Range("I62").Value = 15000
$ws.Range("J62").Value = 37500
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 37500
$ws.Range("M62").Value = -14376
$ws.Range("N62").Value = -38748

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 26250
$ws.Range("I65").Value = 15000
$ws.Range("J65").Value = 37500
$ws.Range("K65").Value = 45000
$ws.Range("L65").Value = 112500
$ws.Range("M65").Value = -41880
$ws.Range("N65").Value = -118740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2941.9092
$ws.Range("I136").Value = 2262.3333
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 6786.999899999999
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -4236.999899999999
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2500000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 2500000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2500000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -2500228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20501.72
$ws.Range("I136").Value = 21297.625
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 63892.875
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -61342.875
$ws.Range("N136").Value = -9300
